# The only substantive (non-cosmetic) change in the target revision is that
# the backlog line "Ajouter un champ statut pour la séance avec une case à
# cocher eleve / encadrant." gets struck through (it was resolved/obsoleted
# by the "evaluation expiration en 7j" work). Everything else in the diff is
# Word's automatic proofing engine re-flagging spelling/grammar ("eleve",
# "seance", "etoile", "evaluation", "mouss", "prenom", "des section", "la
# section choisi", "une popup", "commentaire", "adhérents ,", "email" …) —
# those <w:proofErr/> markers and the run-splits around them carry no text
# or formatting change at all (just cosmetic squiggly-line bookkeeping), so
# we reproduce the one real edit here.

$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "Ajouter un champ statut pour la séance avec une case à cocher eleve / encadrant."
$found = $rng.Find.Execute()

if ($found) {
    # Grab the whole containing paragraph (including its end-of-paragraph
    # mark) so the strikethrough also lands on the paragraph mark run
    # properties, matching how Word records a full-paragraph format change.
    $para = $rng.Paragraphs(1).Range
    $para.Font.StrikeThrough = 1
}
